{"js": "// The author expanded the intro sentence describing what doorPauseControl\n// controls: \"...door and controlling various machine components.\" became\n// \"...door and controlling the Fan motor heater and water inlet valve in\n// the machine.\" Only the FIRST occurrence (end of the opening paragraph)\n// changes; a later, unrelated sentence that also contains the phrase\n// \"various machine components\" (\"It stores the current states of various\n// machine components into variables...\") must stay untouched, so we\n// search for the exact phrase including the trailing period, which is\n// unique to the opening paragraph.\nconst target = \"various machine components.\";\nconst results = context.document.body.search(target, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(`expected exactly 1 match for ${JSON.stringify(target)}, found ${results.items.length}`);\n}\n\nresults.items[0].insertText(\n  \"the Fan motor heater and water inlet valve in the machine.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# The author expanded the intro sentence describing what doorPauseControl\n# controls: \"...door and controlling various machine components.\" became\n# \"...door and controlling the Fan motor heater and water inlet valve in\n# the machine.\" Only the FIRST occurrence (end of the opening paragraph)\n# changes; a later, unrelated sentence elsewhere in the document also\n# contains the phrase \"various machine components\" (\"It stores the\n# current states of various machine components into variables...\") and\n# must stay untouched, so we search for the exact phrase including the\n# trailing period, which is unique to the opening paragraph.\n$d = $word.ActiveDocument\n\n$searchText = \"various machine components.\"\n$replaceText = \"the Fan motor heater and water inlet valve in the machine.\"\n\n# Sanity-check there is exactly one match before touching the document.\n$probe = $d.Content\n$matches = 0\nwhile ($probe.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)) {\n    $matches++\n    $probe.Collapse(0)\n    $probe.SetRange($probe.End, $d.Content.End)\n}\nif ($matches -ne 1) {\n    throw \"expected exactly 1 match for '$searchText', found $matches\"\n}\n\n$rng = $d.Content\n$rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n"}
